$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(13)
$rng = $p.Range
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
$search = " – " + "View" + " Task" + "s"
$found = $rng.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, $search, 2)
Write-Output "found=$found"
